# Update cryptos list cell values (refresh run on Fri Nov 24 10:39:26 UTC 2023 via GitHub Actions)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.896.63"
$ws.Range("E2").Value = "  +1.00%  "
$ws.Range("D3").Value = "2.119.15"
$ws.Range("E3").Value = "  +2.10%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'236.49"
$ws.Range("E5").Value = "  +0.59%  "
$ws.Range("E6").Value = "  +0.73%  "
$ws.Range("D7").Value = "'58.80"
$ws.Range("E7").Value = "  +0.92%  "
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("D9").Value = "'0.393"
$ws.Range("E9").Value = "  +1.89%  "
$ws.Range("D10").Value = "'0.0786"
$ws.Range("E10").Value = "  +2.70%  "
$ws.Range("D12").Value = "2.432.77"
$ws.Range("E12").Value = "  +2.11%  "
$ws.Range("D13").Value = "'14.63"
$ws.Range("E13").Value = "  +1.04%  "
$ws.Range("D14").Value = "'21.43"
$ws.Range("E14").Value = "  +1.75%  "
$ws.Range("D15").Value = "'0.792"
$ws.Range("E15").Value = "  +1.35%  "
$ws.Range("D16").Value = "'5.25"
$ws.Range("E16").Value = "  +1.00%  "
$ws.Range("D17").Value = "2.111.03"
$ws.Range("E17").Value = "  +1.76%  "
$ws.Range("D18").Value = "37.827.94"
$ws.Range("E18").Value = "  +0.29%  "
$ws.Range("D19").Value = "'6.23"
$ws.Range("E19").Value = "  +0.44%  "
$ws.Range("D20").Value = "'70.41"
$ws.Range("E20").Value = "  +1.27%  "
$ws.Range("D21").Value = "0.0₃0827"
$ws.Range("E21").Value = "  +1.15%  "
$ws.Range("D22").Value = "'228.18"
$ws.Range("E22").Value = "  +0.86%  "
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("E24").Value = "  -0.21%  "
$ws.Range("D25").Value = "'2.44"
$ws.Range("E25").Value = "  -2.66%  "
$ws.Range("D26").Value = "'168.78"
$ws.Range("E26").Value = "  +0.99%  "
$ws.Range("D27").Value = "'9.02"
$ws.Range("E27").Value = "  -0.19%  "
$ws.Range("E28").Value = "  +4.67%  "
$ws.Range("E29").Value = "  -4.04%  "
$ws.Range("D30").Value = "'19.52"
$ws.Range("E30").Value = "  +1.42%  "
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("E32").Value = "  +2.37%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.0624"
$ws.Range("E33").Value = "  -0.43%  "
$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").Value = "'2.58"
$ws.Range("E34").Value = "  -0.12%  "
$ws.Range("D35").Value = "'4.60"
$ws.Range("E35").Value = "  +0.15%  "
$ws.Range("E36").Value = "  +5.60%  "
$ws.Range("E37").Value = "  +1.11%  "
$ws.Range("E38").Value = "  -0.23%  "
$ws.Range("D39").Value = "'5.70"
$ws.Range("E39").Value = "  -5.08%  "
$ws.Range("E40").Value = "  -0.10%  "
$ws.Range("D41").Value = "'0.0967"
$ws.Range("E41").Value = "  +1.30%  "
$ws.Range("D42").Value = "'98.09"
$ws.Range("E42").Value = "  +2.59%  "
$ws.Range("D43").Value = "1.479.49"
$ws.Range("E43").Value = "  +0.67%  "
$ws.Range("E44").Value = "  +0.75%  "
$ws.Range("E45").Value = "  -0.51%  "
$ws.Range("E46").Value = "  -8.04%  "
$ws.Range("E47").Value = "  +1.83%  "
$ws.Range("D48").Value = "'15.74"
$ws.Range("E48").Value = "  -0.88%  "
$ws.Range("B49").Value = "MXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D49").Value = "'3.06"
$ws.Range("E49").Value = "  +4.00%  "
$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").Value = "'7.36"
$ws.Range("E50").Value = "  +2.20%  "
$ws.Range("D51").Value = "2.317.02"
$ws.Range("E51").Value = "  +2.12%  "
